# Update status text "Ready for handoff" -> "In Translation" on every
# sheet that shows it (Overview, zh-cn, de-de), and narrow the
# now-shorter status columns to match.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F2").Value = "In Translation"
$overview.Range("E:F").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C:C").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C:C").ColumnWidth = 12.5
